$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; this shifts the existing rows 10..131
# down to 11..132 and carries the row-10 formatting (incl. the date-column
# style) down with them, matching the dimension growing from R131 to R132.
$ws.Rows(10).Insert()

# Populate the newly inserted row 10 with the new price-report record.
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "Macroferia Regional de Talca"
$ws.Range("C10").Value = "Maule"
$ws.Range("D10").Value = 44537
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 100112024
$ws.Range("G10").Value = "Choclo"
$ws.Range("H10").Value = "Choclero"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 20000
$ws.Range("K10").Value = 350
$ws.Range("L10").Value = 350
$ws.Range("M10").Value = 350
$ws.Range("N10").Value = "`$/unidad"
$ws.Range("O10").Value = "Región de O'Higgins"
$ws.Range("P10").Value = 350
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = "Hortaliza"
